# The commit swaps the two embedded themes: the deck's main theme (theme1.xml,
# currently the "Integral" palette used by the slide master) is replaced by the
# stock "Office Theme" color palette (previously only embedded as theme2.xml,
# the notes-master theme). Font scheme / format scheme are identical between
# the two themes already, so only the 12 color-scheme slots actually change.
#
# Apply the new palette via the theme color scheme object model -- this is
# the supported way to rewrite a:clrScheme's srgbClr values through the
# PowerPoint COM surface.

$p  = $ppt.ActivePresentation
$t  = $p.SlideMaster.Theme
$cs = $t.ThemeColorScheme

# Office Theme palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink),
# expressed as the packed 0x00BBGGRR integers the RGB property expects.
$cs.Colors(1).RGB  = 0          # dk1      000000
$cs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$cs.Colors(3).RGB  = 6968388    # dk2      44546A
$cs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$cs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$cs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$cs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$cs.Colors(8).RGB  = 49407      # accent4  FFC000
$cs.Colors(9).RGB  = 12874308   # accent5  4472C4
$cs.Colors(10).RGB = 4697456    # accent6  70AD47
$cs.Colors(11).RGB = 12673797   # hlink    0563C1
$cs.Colors(12).RGB = 7491477    # folHlink 954F72
